# Update numeric values in column F (view/sales counters) across sheets
# to reflect a newer data snapshot, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (sheet "展览")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 199
$ws1.Range("F6").Value  = 580
$ws1.Range("F9").Value  = 495
$ws1.Range("F10").Value = 200
$ws1.Range("F13").Value = 104
$ws1.Range("F14").Value = 520
$ws1.Range("F16").Value = 1656
$ws1.Range("F17").Value = 272
$ws1.Range("F18").Value = 1050
$ws1.Range("F20").Value = 467
$ws1.Range("F21").Value = 13

# 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5206
$ws3.Range("F3").Value = 282

# 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5206
$ws4.Range("F4").Value  = 282
$ws4.Range("F10").Value = 199
$ws4.Range("F17").Value = 580
$ws4.Range("F22").Value = 495
$ws4.Range("F23").Value = 200
$ws4.Range("F27").Value = 104
$ws4.Range("F30").Value = 520
$ws4.Range("F33").Value = 1656
$ws4.Range("F34").Value = 272
$ws4.Range("F35").Value = 1050
$ws4.Range("F38").Value = 467
$ws4.Range("F39").Value = 13
